$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update row labels in column A (respondent-letter suffixes) ---
$ws.Range("A7").Value = "Madre Helena (A)"
$ws.Range("A4").Value = "Papá (B)"
$ws.Range("A5").Value = "Mamá (C)"
$ws.Range("A6").Value = "Hermana (D)"
$ws.Range("A3").Value = "Antonio (E)"

# --- Fill in survey answers B3:K7 (rows 3-7, columns B-K = Pregunta 1..10) ---
$answers3 = @(4, 2, 4, 4, 4, 2, 5, 2, 3, 3)
$answers4 = @(5, 1, 5, 4, 5, 2, 5, 1, 5, 2)
$answers5 = @(4, 1, 4, 3, 4, 2, 4, 2, 3, 2)
$answers6 = @(4, 1, 5, 2, 4, 2, 4, 2, 4, 1)
$answers7 = @(4, 3, 4, 4, 4, 2, 4, 2, 2, 3)

$rowsData = @(
    @{ Row = 3; Values = $answers3 },
    @{ Row = 4; Values = $answers4 },
    @{ Row = 5; Values = $answers5 },
    @{ Row = 6; Values = $answers6 },
    @{ Row = 7; Values = $answers7 }
)

foreach ($entry in $rowsData) {
    $row = $entry.Row
    $vals = $entry.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i  # column B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}

# --- Update Q2 formula (now averages O3:O7 over 5 instead of O2:O7 over 6) ---
$ws.Range("Q2").Formula = "=SUM(O3:O7)/5"

# --- Update selection / active cell ---
$ws.Range("Q3").Select()

# --- Autofit column A to match new bestFit width ---
$ws.Columns.Item(1).AutoFit() | Out-Null
